$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 5.355
$ws.Range("A9").Value = -21.743
$ws.Range("B12").Value = 5.57
$ws.Range("D13").Value = -7.726999999999999
$ws.Range("C15").Value = -13.247
$ws.Range("D16").Value = -8.543000000000001
$ws.Range("A18").Value = -22.051
$ws.Range("A20").Value = -20.793
$ws.Range("D20").Value = -7.558
$ws.Range("D24").Value = -7.547
$ws.Range("B26").Value = 5.505
$ws.Range("A27").Value = -21.4
$ws.Range("B27").Value = 5.770999999999999
$ws.Range("B29").Value = 5.669
$ws.Range("B37").Value = 8.670999999999999
$ws.Range("B38").Value = 5.038
$ws.Range("C38").Value = -12.566
$ws.Range("D39").Value = -7.507
$ws.Range("C44").Value = -12.248
$ws.Range("D48").Value = -7.366000000000001
$ws.Range("B51").Value = 5.790999999999999
$ws.Range("C51").Value = -11.915
$ws.Range("D52").Value = -7.811
$ws.Range("B55").Value = 5.867
$ws.Range("D56").Value = -7.825
$ws.Range("C57").Value = -13.613
$ws.Range("C63").Value = -12.488
$ws.Range("A69").Value = -21.375
$ws.Range("B69").Value = 6.165000000000001
$ws.Range("B70").Value = 5.606
$ws.Range("C70").Value = -11.206
$ws.Range("A76").Value = -20.718
$ws.Range("A82").Value = -22.101
$ws.Range("B83").Value = 6.248
$ws.Range("D84").Value = -8.164
$ws.Range("C99").Value = -12.474
$ws.Range("D100").Value = -8.310999999999998
$ws.Range("D101").Value = -7.831
$ws.Range("B102").Value = 7.398999999999999
